$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ProQuest): hit count & date refreshed ---
$ws.Range("D2").Value = 40
$ws.Range("E2").Value = 41755

# --- Row 3 (EBSCO): link refreshed with new session id; hit count & date refreshed ---
$ws.Range("C3").Value = @'
http://web.a.ebscohost.com/ehost/resultsadvanced?sid=fbdefb89-6f81-40e9-bcec-68c5b3a0f64c%40sessionmgr4004&vid=14&hid=4206
'@
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 41755

# --- Row 4 (ACM): link refreshed (now points at dl.acm.org results); hit count & date refreshed ---
$ws.Range("C4").Value = @'
http://dl.acm.org/results.cfm?within=%28+%0D%0A+Title%3A+%22Augmented+Reality%22+OR+%0D%0A+Abstract%3A+%22Augmented+Reality%22+OR+%0D%0A+Keywords%3A+%22Augmented+Reality%22%0D%0A%29+%0D%0AAND++%0D%0A%28+%0D%0A+%28+%0D%0A++Title%3A+%22College*%22+OR+%0D%0A++Title%3A+%22School*%22+OR+%0D%0A++Title%3A+%22Teach*%22+OR+%0D%0A++Title%3A+%22Learn*%22+OR+%0D%0A++Title%3A+%22Educat*%22+%0D%0A+%29+%0D%0A+OR++%0D%0A+%28+%0D%0A++Abstract%3A+%22College*%22+OR+%0D%0A++Abstract%3A+%22School*%22+OR+%0D%0A++Abstract%3A+%22Teach*%22+OR+%0D%0A++Abstract%3A+%22Learn*%22+OR+%0D%0A++Abstract%3A+%22Educat*%22+%0D%0A+%29+%0D%0A+OR+%0D%0A+%28+%0D%0A++Keywords%3A+%22College*%22+OR+%0D%0A++Keywords%3A+%22School*%22+OR+%0D%0A++Keywords%3A+%22Teach*%22+OR+%0D%0A++Keywords%3A+%22Learn*%22+OR+%0D%0A++Keywords%3A+%22Educat*%22+%0D%0A+%29+%0D%0A%29%0D%0AAND%0D%0A%28%0D%0A%28+%0D%0A++Title%3A+%22Benefi*%22+OR+%0D%0A++Title%3A+%22Advantag*%22++%0D%0A+%29+%0D%0A+OR++%0D%0A+%28+%0D%0A++Abstract%3A+%22Benefi*%22+OR+%0D%0A++Abstract%3A+%22Advantag*%22+%0D%0A+%29+%0D%0A+OR+%0D%0A+%28+%0D%0A++Keywords%3A+%22Benefi*%22+OR+%0D%0A++Keywords%3A+%22Advantag*%22+%0D%0A+%29+%0D%0A%29+&CFID=327269280&CFTOKEN=59077846&adv=1&COLL=DL&qrycnt=405427&DL=ACM&Go.x=48&Go.y=17&termzone=all&allofem=&anyofem=&noneofem=&peoplezone=Name&people=&peoplehow=and&keyword=&keywordhow=AND&affil=&affilhow=AND&pubin=&pubinhow=and&pubby=&pubbyhow=OR&since_year=&before_year=&pubashow=OR&sponsor=&sponsorhow=AND&confdate=&confdatehow=OR&confloc=&conflochow=OR&isbnhow=OR&isbn=&doi=&ccs=&subj=
'@
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 41755

# --- Row 5 (AISeL): query reformatted, working search link added, error comment removed ---
$ws.Range("B5").Value = @'
( 
 ( 
  title:( "Augmented Reality" ) OR 
  abstract:( "Augmented Reality" ) OR 
  subject:( "Augmented Reality" ) 
 ) 
 AND 
 ( 
  title:(
   "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" 
  ) OR 
  abstract:( 
   "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" 
  ) OR 
  subject:( 
   "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" 
  ) 
 ) 
) 
'@
$ws.Range("C5").Value = @'
http://aisel.aisnet.org/do/search/?q=(%20%20%20(%20%20%20%20title%3A(%20%22Augmented%20Reality%22%20)%20OR%20%20%20%20abstract%3A(%20%22Augmented%20Reality%22%20)%20OR%20%20%20%20subject%3A(%20%22Augmented%20Reality%22%20)%20%20%20)%20%20%20AND%20%20%20(%20%20%20%20title%3A(%20%20%20%20%22College*%22%20OR%20%22School*%22%20OR%20%22Teach*%22%20OR%20%22Learn*%22%20OR%20%22Educat*%22%20%20%20%20)%20OR%20%20%20%20abstract%3A(%20%20%20%20%20%22College*%22%20OR%20%22School*%22%20OR%20%22Teach*%22%20OR%20%22Learn*%22%20OR%20%22Educat*%22%20%20%20%20)%20OR%20%20%20%20subject%3A(%20%20%20%20%22College*%22%20OR%20%22School*%22%20OR%20%22Teach*%22%20OR%20%22Learn*%22%20OR%20%22Educat*%22%20%20%20%20)%20%20%20)%20%20)%20&start=0&context=509156
'@
$ws.Range("F5").ClearContents()
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 41755

# --- Row 2 (ProQuest): query refreshed last (now also searches Benefit/Advantage terms) ---
$ws.Range("C2").Value = @'
http://search.proquest.com/results/1EE1D94F9BA42E4PQ/1/$5bqueryType$3dcommandline:OS$3b+sortType$3drelevance$3b+searchTerms$3d$5b$3cAND$7call:$28TI$28$22Augmented+Reality$22$29+OR+AB$28$22Augmented+Reality$22$29+OR+SU$28$22Augmented+Reality$22$29$29+AND+$28TI$28$22College*$22+OR+$22School*$22+OR+$22Teach*$22+OR+$22Learn*$22+OR+$22Educat*$22$29+OR+AB$28$22College*$22+OR+$22School*$22+OR+$22Teach*$22+OR+$22Learn*$22+OR+$22Educat*$22$29+OR+SU$28$22College*$22+OR+$22School*$22+OR+$22Teach*$22+OR+$22Learn*$22+OR+$22Educat*$22$29$29+AND+$28TI$28$22Benefi*$22+OR+$22Advantag*$22$29+OR+AB$28$22Benefi*$22+OR+$22Advantag*$22$29+OR+SU$28$22Benefi*$22+OR+$22Advantag*$22$29$29$3e$5d$3b+searchParameters$3d$7bNAVIGATORS$3dnavsummarynav,sourcetypenav,pubtitlenav,objecttypenav,languagenav$28filter$3d200$2f0$2f*$29,decadenav$28filter$3d110$2f0$2f*,sort$3dname$2fascending$29,yearnav$28filter$3d1100$2f0$2f*,sort$3dname$2fascending$29,yearmonthnav$28filter$3d120$2f0$2f*,sort$3dname$2fascending$29,monthnav$28sort$3dname$2fascending$29,daynav$28sort$3dname$2fascending$29,+RS$3dOP,+chunkSize$3d20,+instance$3dprod.academic,+ftblock$3d740842+1+199113+113+670831+670829+660845+199000+660843+199001+660840,+removeDuplicates$3dtrue$7d$3b+metaData$3d$7bUsageSearchMode$3dCommandLine,+dbselections$3dallAvailable,+SEARCH_ID_TIMESTAMP$3d1398524529400$7d$5d?accountid=10218
'@
$ws.Range("B2").Value = @'
 ( TI( "Augmented Reality" ) OR AB( "Augmented Reality" ) OR SU( "Augmented Reality" ) ) AND ( TI( "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" ) OR AB( "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" ) OR SU( "College*" OR "School*" OR "Teach*" OR "Learn*" OR "Educat*" ) ) 
'@
